$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 9 first: a lone space in B9 (this becomes shared-string #12).
$ws.Range("B9").Value = " "

# Row 8: fill in the "comment" column (B8, becomes shared-string #13) and
# score (C8) for TP4, matching the style already used on row 6 (wrap-text /
# filled comment cell).
$ws.Range("B8").Value = "Identificar erros no template method`n- python + antlr`n- boa explicação, custei para entender um pouco o exemplo.`n- Entendi o que o grupo propôs a validar`n- Fez com listener.`n- apontou limitações e dificuldades.`n- saiu da zona de conforto, aprendeu um novo padrão e mandaram ver."
$ws.Range("B6").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("C8").Value = 10

# Row heights / selection tweaks that accompanied the content change.
$ws.Rows.Item(8).RowHeight = 148
$ws.Range("F7").Select()
